$d = $word.ActiveDocument

$d.Content.Find.Execute("812÷4=203, 0", $true, $false, $false, $false, $false, $true, 1, $false, "535÷3=178, 1", 2) | Out-Null
$d.Content.Find.Execute("681÷9=75, 6", $true, $false, $false, $false, $false, $true, 1, $false, "364÷2=182, 0", 2) | Out-Null
$d.Content.Find.Execute("123÷5=24, 3", $true, $false, $false, $false, $false, $true, 1, $false, "564÷9=62, 6", 2) | Out-Null
$d.Content.Find.Execute("896÷6=149, 2", $true, $false, $false, $false, $false, $true, 1, $false, "731÷3=243, 2", 2) | Out-Null
$d.Content.Find.Execute("156÷6=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "289÷2=144, 1", 2) | Out-Null
$d.Content.Find.Execute("716÷5=143, 1", $true, $false, $false, $false, $false, $true, 1, $false, "936÷5=187, 1", 2) | Out-Null
$d.Content.Find.Execute("217÷9=24, 1", $true, $false, $false, $false, $false, $true, 1, $false, "973÷3=324, 1", 2) | Out-Null
$d.Content.Find.Execute("285÷3=95, 0", $true, $false, $false, $false, $false, $true, 1, $false, "485÷2=242, 1", 2) | Out-Null
$d.Content.Find.Execute("841÷2=420, 1", $true, $false, $false, $false, $false, $true, 1, $false, "840÷9=93, 3", 2) | Out-Null
$d.Content.Find.Execute("437÷4=109, 1", $true, $false, $false, $false, $false, $true, 1, $false, "427÷8=53, 3", 2) | Out-Null
$d.Content.Find.Execute("418÷5=83, 3", $true, $false, $false, $false, $false, $true, 1, $false, "282÷6=47, 0", 2) | Out-Null
$d.Content.Find.Execute("533÷9=59, 2", $true, $false, $false, $false, $false, $true, 1, $false, "135÷3=45, 0", 2) | Out-Null
$d.Content.Find.Execute("832÷9=92, 4", $true, $false, $false, $false, $false, $true, 1, $false, "190÷8=23, 6", 2) | Out-Null
$d.Content.Find.Execute("832÷4=208, 0", $true, $false, $false, $false, $false, $true, 1, $false, "852÷8=106, 4", 2) | Out-Null
$d.Content.Find.Execute("341÷8=42, 5", $true, $false, $false, $false, $false, $true, 1, $false, "440÷3=146, 2", 2) | Out-Null
$d.Content.Find.Execute("198÷3=66, 0", $true, $false, $false, $false, $false, $true, 1, $false, "813÷9=90, 3", 2) | Out-Null
$d.Content.Find.Execute("410÷2=205, 0", $true, $false, $false, $false, $false, $true, 1, $false, "425÷4=106, 1", 2) | Out-Null
$d.Content.Find.Execute("288÷9=32, 0", $true, $false, $false, $false, $false, $true, 1, $false, "803÷9=89, 2", 2) | Out-Null
$d.Content.Find.Execute("281÷6=46, 5", $true, $false, $false, $false, $false, $true, 1, $false, "842÷2=421, 0", 2) | Out-Null
$d.Content.Find.Execute("484÷8=60, 4", $true, $false, $false, $false, $false, $true, 1, $false, "586÷7=83, 5", 2) | Out-Null
$d.Content.Find.Execute("847÷4=211, 3", $true, $false, $false, $false, $false, $true, 1, $false, "647÷5=129, 2", 2) | Out-Null
$d.Content.Find.Execute("941÷6=156, 5", $true, $false, $false, $false, $false, $true, 1, $false, "482÷8=60, 2", 2) | Out-Null
$d.Content.Find.Execute("725÷5=145, 0", $true, $false, $false, $false, $false, $true, 1, $false, "747÷5=149, 2", 2) | Out-Null
$d.Content.Find.Execute("740÷9=82, 2", $true, $false, $false, $false, $false, $true, 1, $false, "648÷7=92, 4", 2) | Out-Null
$d.Content.Find.Execute("736÷6=122, 4", $true, $false, $false, $false, $false, $true, 1, $false, "310÷5=62, 0", 2) | Out-Null

Write-Host "Done applying replacements"
